$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.481.42'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.568.47'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.38%  '
$ws.Range("E6").Value = '  -0.67%  '
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.98'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0868'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.65%  '
$ws.Range("D12").Value = '1.791.35'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").Value = '1.571.83'
$ws.Range("E13").Value = '  -0.40%  '
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.518'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.83%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.40'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.43%  '
$ws.Range("D17").Value = '27.459.00'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.20'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.03'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("B26").Value = 'BinanceUSD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.27%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.08%  '
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("E29").Value = '  -1.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.42%  '
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("D33").Value = '1.363.43'
$ws.Range("E33").Value = '  +0.55%  '
$ws.Range("E34").Value = '  +0.93%  '
$ws.Range("E35").Value = '  +2.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.978'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.92%  '
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("E38").Value = '  +2.12%  '
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("E40").Value = '  +2.07%  '
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.974'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  +2.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.76%  '
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("D47").Value = '1.703.63'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.45'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.68%  '
$ws.Range("E49").Value = '  -1.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0954'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.73%  '
$ws.Range("E51").Value = '  -0.16%  '
